# Created FrameWork and Added TestCases for REPORT Module
$wb = $excel.ActiveWorkbook

# The workbook used to hold two login-style test-case sheets (TC001/TC002).
# The new "REPORT module" test framework only needs a single sheet, so drop
# the second one and repurpose/rename the first.
[void]$wb.Worksheets.Item("TC002").Delete()

$ws = $wb.Worksheets.Item("TC001")
$ws.Name = "Sheet1"

# Clear out the old username/password test data (columns B:C) without
# disturbing column A's existing formatting.
[void]$ws.Range("B1:C2").ClearContents()

# New report-chart test data.
$ws.Range("A1").Value = "ChartName"
$ws.Range("A3").Value = "praba"
$ws.Range("A2").Value = "Timing Chart"

[void]$ws.Columns("A:A").AutoFit()

[void]$ws.Range("B5").Select()
